# Update CDA Logical model for ST.r2b
# - Bump Version / Date metadata values
# - Insert a new "Jurisdiction" property row (empty value) between
#   "Contact" and "Description" on the Metadata sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3: Property="Version")
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8: Property="Date")
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row before row 11 ("Description") for "Jurisdiction"
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Match the formatting of the surrounding data rows (Insert() alone
# leaves the new row with a default style instead of the shared
# data-row style used throughout the table)
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
